# Update sheet "data" (sheet1): add column AJ with data for "12. 10. 2021"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# Copy header cell formatting (border/font/alignment) from AI1 into the new AJ1, then set its value/text
$ws1.Range("AI1").Copy($ws1.Range("AJ1"))
$ws1.Range("AJ1").Value = "12. 10. 2021"

# New data column AJ (rows 2-58), values for the "12. 10. 2021" survey wave
$ws1.Range("AJ2").Value = 0.63
$ws1.Range("AJ3").Value = 0.24
$ws1.Range("AJ4").Value = 0.13
$ws1.Range("AJ5").Value = 0.76
$ws1.Range("AJ6").Value = 0.14
$ws1.Range("AJ7").Value = 0.1
$ws1.Range("AJ8").Value = 0.71
$ws1.Range("AJ9").Value = 0.2
$ws1.Range("AJ10").Value = 0.09
$ws1.Range("AJ11").Value = 0.47
$ws1.Range("AJ12").Value = 0.35
$ws1.Range("AJ13").Value = 0.18
$ws1.Range("AJ14").Value = 0.45
$ws1.Range("AJ15").Value = 0.33
$ws1.Range("AJ16").Value = 0.22
$ws1.Range("AJ17").Value = 0.5600000000000001
$ws1.Range("AJ18").Value = 0.29
$ws1.Range("AJ19").Value = 0.15
$ws1.Range("AJ20").Value = 0.7
$ws1.Range("AJ21").Value = 0.21
$ws1.Range("AJ22").Value = 0.09
$ws1.Range("AJ23").Value = 0.79
$ws1.Range("AJ24").Value = 0.13
$ws1.Range("AJ25").Value = 0.08
$ws1.Range("AJ26").Value = 0.64
$ws1.Range("AJ27").Value = 0.25
$ws1.Range("AJ28").Value = 0.11
$ws1.Range("AJ29").Value = 0.45
$ws1.Range("AJ30").Value = 0.36
$ws1.Range("AJ31").Value = 0.19
$ws1.Range("AJ32").Value = 0.62
$ws1.Range("AJ33").Value = 0.25
$ws1.Range("AJ34").Value = 0.13
$ws1.Range("AJ35").Value = 0.5600000000000001
$ws1.Range("AJ36").Value = 0.27
$ws1.Range("AJ37").Value = 0.17
$ws1.Range("AJ38").Value = 0.65
$ws1.Range("AJ39").Value = 0.23
$ws1.Range("AJ40").Value = 0.12
$ws1.Range("AJ41").Value = 0.65
$ws1.Range("AJ42").Value = 0.22
$ws1.Range("AJ43").Value = 0.13
$ws1.Range("AJ44").Value = 0.61
$ws1.Range("AJ45").Value = 0.26
$ws1.Range("AJ46").Value = 0.13
$ws1.Range("AJ47").Value = 0.67
$ws1.Range("AJ48").Value = 0.2
$ws1.Range("AJ49").Value = 0.13
$ws1.Range("AJ50").Value = 0.59
$ws1.Range("AJ51").Value = 0.29
$ws1.Range("AJ52").Value = 0.12
$ws1.Range("AJ53").Value = 0.61
$ws1.Range("AJ54").Value = 0.27
$ws1.Range("AJ55").Value = 0.12
$ws1.Range("AJ56").Value = 0.5600000000000001
$ws1.Range("AJ57").Value = 0.29
$ws1.Range("AJ58").Value = 0.15

# Update the footnote/caption row with the new "aktualizace" date
$ws1.Range("A59").Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"


# Update sheet "pocetR" (sheet2): add column AI with sample-size data for "12. 10. 2021"
$ws2 = $wb.Worksheets.Item("pocetR")

# Copy header cell formatting from AH1 into the new AI1, then set its value/text
$ws2.Range("AH1").Copy($ws2.Range("AI1"))
$ws2.Range("AI1").Value = "12. 10. 2021"

# New data column AI (rows 2-20), sample sizes for the "12. 10. 2021" survey wave
$ws2.Range("AI2").Value = 1836
$ws2.Range("AI3").Value = 454
$ws2.Range("AI4").Value = 670
$ws2.Range("AI5").Value = 712
$ws2.Range("AI6").Value = 272
$ws2.Range("AI7").Value = 490
$ws2.Range("AI8").Value = 812
$ws2.Range("AI9").Value = 522
$ws2.Range("AI10").Value = 819
$ws2.Range("AI11").Value = 495
$ws2.Range("AI12").Value = 296
$ws2.Range("AI13").Value = 316
$ws2.Range("AI14").Value = 1224
$ws2.Range("AI15").Value = 895
$ws2.Range("AI16").Value = 941
$ws2.Range("AI17").Value = 960
$ws2.Range("AI18").Value = 418
$ws2.Range("AI19").Value = 216
$ws2.Range("AI20").Value = 242

# Update the footnote/caption row with the new "aktualizace" date
$ws2.Range("A21").Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"
# Row 21 also carries a trailing empty-string cell pattern across the table; extend it to the new AI column
$ws2.Range("AI21").Value = ""

